$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New weekly data point: insert two new rows at the top of this item's
# price block (rows 80-81), pushing the existing rows (old 80-147) down
# to 82-149.
$ws.Rows("80:81").Insert()

# Row 80 - Ají "Americana (o)" / "Primera"
$ws.Range("A80").Value = 2
$ws.Range("B80").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C80").Value = "Coquimbo"
$ws.Range("D80").Value = 44447
$ws.Range("E80").Value = 4
$ws.Range("F80").Value = 100112021
$ws.Range("G80").Value = "Ají"
$ws.Range("H80").Value = "Americana (o)"
$ws.Range("I80").Value = "Primera"
$ws.Range("J80").Value = 140
$ws.Range("K80").Value = 65000
$ws.Range("L80").Value = 70000
$ws.Range("M80").Value = 67500
$ws.Range("N80").Value = "$/caja 25 kilos"
$ws.Range("O80").Value = "Provincia de Limarí"
$ws.Range("P80").Value = 2700
$ws.Range("Q80").Value = 25
$ws.Range("R80").Value = "Hortaliza"

# Row 81 - Ají "Americana (o)" / "Segunda"
$ws.Range("A81").Value = 2
$ws.Range("B81").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C81").Value = "Coquimbo"
$ws.Range("D81").Value = 44447
$ws.Range("E81").Value = 4
$ws.Range("F81").Value = 100112021
$ws.Range("G81").Value = "Ají"
$ws.Range("H81").Value = "Americana (o)"
$ws.Range("I81").Value = "Segunda"
$ws.Range("J81").Value = 100
$ws.Range("K81").Value = 55000
$ws.Range("L81").Value = 60000
$ws.Range("M81").Value = 57500
$ws.Range("N81").Value = "$/caja 25 kilos"
$ws.Range("O81").Value = "Provincia de Limarí"
$ws.Range("P81").Value = 2300
$ws.Range("Q81").Value = 25
$ws.Range("R81").Value = "Hortaliza"
